$wb = $excel.ActiveWorkbook

# Sheet1: update Total Stories (F2) and Open Defects (K2) for the first project row
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("F2").Value = 4
$ws1.Range("K2").Value = 8

# Text_Summary_REPORT: update the matching summary text lines
$ws2 = $wb.Worksheets.Item("Text_Summary_REPORT")
$ws2.Range("A5").Value = "4 {'story' if total_stories == 1 else 'stories'} in scope + Regression."
$ws2.Range("A6").Value = "* 8 defect(s) still open and 0 defect(s) closed today."
